$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Update the "troposphere tau" sheet (sheet1) flight labels, which
# were re-sorted after recomputing the tropospheric lifetime values ---

$ws1.Range("F4").Value2 = "RF06"
$ws1.Range("F5").Value2 = "RF07"
$ws1.Range("F6").Value2 = "RF05"
$ws1.Range("F7").Value2 = "RF09"
$ws1.Range("F8").Value2 = "RF10"
$ws1.Range("F9").Value2 = "RF08"

$ws1.Range("F10").Value2 = "RF13"
$ws1.Range("H10").Value2 = 14
$ws1.Range("I10").Value2 = 0.92

$ws1.Range("F11").Value2 = "RF11"
$ws1.Range("H11").Value2 = 14.2
$ws1.Range("I11").Value2 = 0.9

$ws1.Range("F12").Value2 = "RF12"

# New blank separator cell at J13 matching the style used by the rest of
# column E/J (the thin blank divider column between the two tables).
$ws1.Range("E4").Copy()
$ws1.Range("J13").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Update stored selections on both sheets (order matters: select the
# inactive sheet first so that sheet1 remains the active tab) ---
[void]$ws2.Range("E16").Select()
[void]$ws1.Range("E15").Select()
